# Delete description from each excel sheet
# (Title / Update Date / Comment header rows 1-3) on the template sheets,
# leaving the "Note" sheet's content untouched.

$wb = $excel.ActiveWorkbook

# --- access-template: remove the 3 description rows (rows 1-3) ---
$ws2 = $wb.Worksheets.Item("access-template")
[void]$ws2.Rows("1:3").Delete()
[void]$ws2.Range("A1:XFD3").Select()

# --- snmp-template: remove the 3 description rows (rows 1-3) ---
$ws3 = $wb.Worksheets.Item("snmp-template")
[void]$ws3.Rows("1:3").Delete()
[void]$ws3.Range("A1").Select()

# --- Note: selection simply moved one column to the right, and the
#     sheet now carries an explicit (portrait / letter-ish) page setup ---
$ws1 = $wb.Worksheets.Item("Note")
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
$ws1.Activate()
[void]$ws1.Range("C6").Select()
